$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up existing rows (random index increments shifted a couple of values) ---
# Row 12 (Person_381 / Acro): Ropes flips from TRUE to FALSE
$ws.Range("E12").Value = $false

# Row 37 (Person_792 / Legends ): Tag gains a "Leadership" value
$ws.Range("F37").Value = "Leadership"

# --- Append 4 new staff rows (43-46) ---
$ws.Range("A43").Value = "Person_798"
$ws.Range("B43").Value = "Socks"
$ws.Range("C43").Value = "Male"
$ws.Range("D43").Value = $false
$ws.Range("E43").Value = $false
$ws.Range("F43").Value = "Leadership"

$ws.Range("A44").Value = "Person_799"
$ws.Range("B44").Value = '"2319"'
$ws.Range("C44").Value = "Male"
$ws.Range("D44").Value = $false
$ws.Range("E44").Value = $false
$ws.Range("F44").Value = "Leadership"

$ws.Range("A45").Value = "Person_800"
$ws.Range("B45").Value = "Swamp Puppy"
$ws.Range("C45").Value = "Male"
$ws.Range("D45").Value = $false
$ws.Range("E45").Value = $false
$ws.Range("F45").Value = "Leadership"

$ws.Range("A46").Value = "Person_801"
$ws.Range("B46").Value = "Hullabaloo"
$ws.Range("C46").Value = "Female"
$ws.Range("D46").Value = $false
$ws.Range("E46").Value = $false
$ws.Range("F46").Value = "Leadership"
